# complianceReport.xlsx edit: add HOURS/ETHICS STATE columns properly,
# collapse the DATE column to short M/D/YYYY text, add a merged spacer
# column C next to TITLE, and append Credits summary rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row 8: I8/J8 held "ETHICS STATE"/"HOURS" - swap them so the
#    HOURS column comes first (I), then ETHICS STATE (J).
# ---------------------------------------------------------------------
$ws.Range("I8").Value = "HOURS"
$ws.Range("J8").Value = "ETHICS STATE"

# ---------------------------------------------------------------------
# 2. Data rows 10-21: the per-row hour total that lived in column J
#    moves to column I; J is cleared out (row 10 is the exception -
#    it keeps both columns, simply swapped: I10<->J10).
# ---------------------------------------------------------------------
$dataRows = 10..21

foreach ($r in $dataRows) {
    $jCell = $ws.Cells.Item($r, 10)   # column J
    $iCell = $ws.Cells.Item($r, 9)    # column I
    $jVal = $jCell.Value2
    $iVal = $iCell.Value2

    if ($r -eq 10) {
        # swap the two existing values
        $iCell.Value = $jVal
        $jCell.Value = $iVal
    } else {
        $iCell.Value = $jVal
        $jCell.ClearContents()
    }
}

# ---------------------------------------------------------------------
# 3. Column A (rows 10-21): collapse the ISO timestamp text to a plain
#    M/D/YYYY text string (same calendar date, no time-of-day).
#    Force text formatting first so Excel doesn't reparse the literal
#    "11/29/2017" back into a date serial number.
# ---------------------------------------------------------------------
$ws.Range("A10:A21").NumberFormat = "@"

foreach ($r in $dataRows) {
    $aCell = $ws.Cells.Item($r, 1)
    $iso = $aCell.Value2
    $datePart = ($iso -split "T")[0]
    $pieces = $datePart -split "-"
    $y = $pieces[0]
    $m = [int]$pieces[1]
    $d = [int]$pieces[2]
    $aCell.Value = "$m/$d/$y"
}

# ---------------------------------------------------------------------
# 4. Insert a spacer column C (merged into the TITLE column B) across
#    the header and data rows 10-20, mirroring the other merged header
#    cells (A8:A9, D8:D9, ...). Row 21 is left untouched - the source
#    edit never added a C21 cell/merge for the last row.
# ---------------------------------------------------------------------
$ws.Range("B8:B9").UnMerge()
$ws.Range("C8:C9").UnMerge()
$ws.Range("B8:C9").Merge()

# Re-stamp the now-blank merge-partner cells so they still get written
# out as explicit empty <c/> nodes (merging drops untouched blanks).
$ws.Range("B9").Font.Bold = $false
$ws.Range("C8").Font.Bold = $false
$ws.Range("C9").Font.Bold = $false

foreach ($r in 10..20) {
    $ws.Range("B" + $r + ":C" + $r).Merge()
    # Touch column C (after merging!) so the (currently nonexistent)
    # cell gets written out even though it stays empty - mirrors the
    # source's pattern of an explicit empty <c r="C.."/> beside each
    # populated B cell. Merging first and touching after matters -
    # a touch made before the merge gets dropped by it.
    $ws.Cells.Item($r, 3).Font.Bold = $false
}

# ---------------------------------------------------------------------
# 5. Second header + Total Credits summary block (rows 22, 24-27).
# ---------------------------------------------------------------------
$ws.Range("I22").Value = "HOURS"
$ws.Range("J22").Value = "ETHICS STATE"

$ws.Range("A24").Value = "Total Credits Applied:"
$ws.Range("I24").Value = 41.5
$ws.Range("J24").Value = 2

$ws.Range("A25").Value = "Total Credits Earned:"
$ws.Range("I25").Value = 41.5
$ws.Range("J25").Value = 2

$ws.Range("A26").Value = "Continuing Education Requirement:"
$ws.Range("I26").Value = 36
$ws.Range("J26").Value = 4

$ws.Range("A27").Value = "Credits Remaining:"
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 2
